$wb = $excel.ActiveWorkbook

function Set-Cell($sheetObj, $addr, $val) {
    if ($null -eq $val) {
        $sheetObj.Range($addr).ClearContents()
    } else {
        $sheetObj.Range($addr).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws "H21" 7999.75
Set-Cell $ws "I21" 3999.6667
Set-Cell $ws "J21" 20000
Set-Cell $ws "K21" 3999.6667
Set-Cell $ws "L21" 20000
Set-Cell $ws "M21" -3531.6667
Set-Cell $ws "N21" -20936
Set-Cell $ws "H23" 7999.75
Set-Cell $ws "I23" 3999.6667
Set-Cell $ws "J23" 20000
Set-Cell $ws "K23" 3999.6667
Set-Cell $ws "L23" 20000
Set-Cell $ws "M23" -3765.6667
Set-Cell $ws "N23" -20468
Set-Cell $ws "H40" 304698.8
Set-Cell $ws "I40" 2998
Set-Cell $ws "J40" 606399.6
Set-Cell $ws "K40" 2998
Set-Cell $ws "L40" 606399.6
Set-Cell $ws "M40" -2823
Set-Cell $ws "N40" -606749.6
Set-Cell $ws "H98" 0
Set-Cell $ws "I98" 0
Set-Cell $ws "J98" 0
Set-Cell $ws "K98" 0
Set-Cell $ws "L98" 0
Set-Cell $ws "M98" $null
Set-Cell $ws "N98" $null
Set-Cell $ws "H115" 1000
Set-Cell $ws "I115" 1000
Set-Cell $ws "K115" 3000
Set-Cell $ws "M115" -1433
Set-Cell $ws "H116" 9000
Set-Cell $ws "I116" 9000
Set-Cell $ws "J116" 0
Set-Cell $ws "K116" 9000
Set-Cell $ws "L116" 0
Set-Cell $ws "M116" -5558
Set-Cell $ws "N116" $null
Set-Cell $ws "H122" 0
Set-Cell $ws "I122" 0
Set-Cell $ws "J122" 0
Set-Cell $ws "K122" 0
Set-Cell $ws "L122" 0
Set-Cell $ws "M122" $null
Set-Cell $ws "N122" $null
Set-Cell $ws "H129" 2599.4285
Set-Cell $ws "J129" 1500
Set-Cell $ws "L129" 4500
Set-Cell $ws "N129" -14500
Set-Cell $ws "H137" 2866
Set-Cell $ws "I137" 2866
Set-Cell $ws "K137" 8598
Set-Cell $ws "M137" -6048
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws "H32" 6123.756
Set-Cell $ws "I32" 6270.3945
Set-Cell $ws "K32" 6270.3945
Set-Cell $ws "M32" -5983.3945
Set-Cell $ws "H61" 4222.8
Set-Cell $ws "I61" 4222.8
Set-Cell $ws "K61" 4222.8
Set-Cell $ws "M61" -4010.8
Set-Cell $ws "H74" 7610.0713
Set-Cell $ws "I74" 7651.478
Set-Cell $ws "J74" 7419.6
Set-Cell $ws "K74" 7651.478
Set-Cell $ws "L74" 7419.6
Set-Cell $ws "M74" -6777.478
Set-Cell $ws "N74" -9167.6
Set-Cell $ws "H77" 7610.0713
Set-Cell $ws "I77" 7651.478
Set-Cell $ws "J77" 7419.6
Set-Cell $ws "K77" 38257.39
Set-Cell $ws "L77" 37098
Set-Cell $ws "M77" -33889.39
Set-Cell $ws "N77" -45834
Set-Cell $ws "H97" 3999.3333
Set-Cell $ws "I97" 3999.5
Set-Cell $ws "J97" 3999
Set-Cell $ws "K97" 3999.5
Set-Cell $ws "L97" 3999
Set-Cell $ws "M97" -3503.5
Set-Cell $ws "N97" -4991
Set-Cell $ws "H122" 2742
Set-Cell $ws "I122" 2742
Set-Cell $ws "J122" 0
Set-Cell $ws "K122" 8226
Set-Cell $ws "L122" 0
Set-Cell $ws "M122" -5776
Set-Cell $ws "N122" $null
Set-Cell $ws "H136" 4222.8
Set-Cell $ws "I136" 4222.8
Set-Cell $ws "K136" 12668.4
Set-Cell $ws "M136" -10118.4
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws "H20" 11276.75
Set-Cell $ws "J20" 20499.5
Set-Cell $ws "L20" 20499.5
Set-Cell $ws "N20" -20993.5
Set-Cell $ws "H22" 639.8570999999999
Set-Cell $ws "I22" 500
Set-Cell $ws "K22" 500
Set-Cell $ws "M22" -327
Set-Cell $ws "H94" 3874
Set-Cell $ws "I94" 2284.5715
Set-Cell $ws "J94" 15000
Set-Cell $ws "K94" 2284.5715
Set-Cell $ws "L94" 15000
Set-Cell $ws "M94" -1833.5715
Set-Cell $ws "N94" -15902
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws "H132" 3247.3635
Set-Cell $ws "I132" 2602.875
Set-Cell $ws "K132" 7808.625
Set-Cell $ws "M132" -5278.625
Set-Cell $ws "H134" 7329.2856
Set-Cell $ws "I134" 7329.2856
Set-Cell $ws "K134" 21987.8568
Set-Cell $ws "M134" -19452.8568
$ws = $wb.Worksheets.Item("CUL")
Set-Cell $ws "H114" 1264.8334
Set-Cell $ws "I114" 1566.3334
Set-Cell $ws "J114" 963.3333
Set-Cell $ws "K114" 4699.0002
Set-Cell $ws "L114" 2889.9999
Set-Cell $ws "M114" -1445.0002
Set-Cell $ws "N114" -9397.999899999999
Set-Cell $ws "H117" 914
Set-Cell $ws "I117" 495.33334
Set-Cell $ws "J117" 1228
Set-Cell $ws "K117" 1486.00002
Set-Cell $ws "L117" 3684
Set-Cell $ws "M117" 1955.99998
Set-Cell $ws "N117" -10568
Set-Cell $ws "H121" 137.8
Set-Cell $ws "I121" 137.8
Set-Cell $ws "K121" 413.4
Set-Cell $ws "M121" 896.5999999999999
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws "H132" 3698.7778
Set-Cell $ws "I132" 2698.2
Set-Cell $ws "J132" 4949.5
Set-Cell $ws "K132" 8094.599999999999
Set-Cell $ws "L132" 14848.5
Set-Cell $ws "M132" -5564.599999999999
Set-Cell $ws "N132" -19908.5
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws "H16" 1900
Set-Cell $ws "I16" 1900
Set-Cell $ws "K16" 1900
Set-Cell $ws "M16" -1730
Set-Cell $ws "H100" 5560.8
Set-Cell $ws "I100" 6350
Set-Cell $ws "J100" 2404
Set-Cell $ws "K100" 6350
Set-Cell $ws "L100" 2404
Set-Cell $ws "M100" -5809
Set-Cell $ws "N100" -3486
$ws = $wb.Worksheets.Item("WVR")
Set-Cell $ws "H55" 13903
Set-Cell $ws "I55" 0
Set-Cell $ws "J55" 13903
Set-Cell $ws "K55" 0
Set-Cell $ws "L55" 13903
Set-Cell $ws "M55" $null
Set-Cell $ws "N55" -14457
Set-Cell $ws "H74" 59994.5
Set-Cell $ws "I74" 44999
Set-Cell $ws "J74" 64993
Set-Cell $ws "K74" 44999
Set-Cell $ws "L74" 64993
Set-Cell $ws "M74" -44063
Set-Cell $ws "N74" -66865
Set-Cell $ws "H77" 59994.5
Set-Cell $ws "I77" 44999
Set-Cell $ws "J77" 64993
Set-Cell $ws "K77" 134997
Set-Cell $ws "L77" 194979
Set-Cell $ws "M77" -130317
Set-Cell $ws "N77" -204339
Set-Cell $ws "H96" 6959.8
Set-Cell $ws "I96" 6959.8
Set-Cell $ws "K96" 6959.8
Set-Cell $ws "M96" -5586.8
Set-Cell $ws "H122" 4235.1
Set-Cell $ws "I122" 4235.1
Set-Cell $ws "K122" 12705.3
Set-Cell $ws "M122" -10255.3
Set-Cell $ws "H132" 3949.2
Set-Cell $ws "I132" 550
Set-Cell $ws "K132" 1650
Set-Cell $ws "M132" 880
Set-Cell $ws "H136" 3390.8948
Set-Cell $ws "I136" 3384.8333
Set-Cell $ws "J136" 3500
Set-Cell $ws "K136" 10154.4999
Set-Cell $ws "L136" 10500
Set-Cell $ws "M136" -7604.499899999999
Set-Cell $ws "N136" -15600
